# Scheduled-runner refresh: update cached market-price/profit figures
# (columns H:N) for a handful of Leve rows across several crafting-job
# sheets. Source columns A:G are untouched.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2459.1785
$ws.Range("I19").Value = 7511.4287
$ws.Range("J19").Value = 775.0952
$ws.Range("K19").Value = 7511.4287
$ws.Range("L19").Value = 775.0952
$ws.Range("M19").Value = -7336.4287
$ws.Range("N19").Value = -1125.0952

$ws.Range("H112").Value = 1415.1428
$ws.Range("I112").Value = 665
$ws.Range("J112").Value = 1570.3448
$ws.Range("K112").Value = 1995
$ws.Range("L112").Value = 4711.0344
$ws.Range("M112").Value = -887
$ws.Range("N112").Value = -6927.0344

$ws.Range("H116").Value = 2882.4285
$ws.Range("J116").Value = 2963
$ws.Range("L116").Value = 2963
$ws.Range("N116").Value = -9847

$ws.Range("H127").Value = 1003
$ws.Range("I127").Value = 600
$ws.Range("J127").Value = 1189
$ws.Range("K127").Value = 1800
$ws.Range("L127").Value = 3567
$ws.Range("M127").Value = 3160
$ws.Range("N127").Value = -13487

$ws.Range("H132").Value = 4290.6055
$ws.Range("I132").Value = 2118.724
$ws.Range("J132").Value = 11288.889
$ws.Range("K132").Value = 6356.172
$ws.Range("L132").Value = 33866.667
$ws.Range("M132").Value = -3826.172
$ws.Range("N132").Value = -38926.667

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3351.2856
$ws.Range("I63").Value = 3493.1667
$ws.Range("K63").Value = 3493.1667
$ws.Range("M63").Value = -2807.1667

$ws.Range("H66").Value = 3351.2856
$ws.Range("I66").Value = 3493.1667
$ws.Range("K66").Value = 17465.8335
$ws.Range("M66").Value = -14033.8335

$ws.Range("H132").Value = 28082.1
$ws.Range("I132").Value = 35109.547
$ws.Range("J132").Value = 3876.4443
$ws.Range("K132").Value = 105328.641
$ws.Range("L132").Value = 11629.3329
$ws.Range("M132").Value = -102798.641
$ws.Range("N132").Value = -16689.3329

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 40780
$ws.Range("J132").Value = 40780
$ws.Range("L132").Value = 40780
$ws.Range("N132").Value = -50900

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2046.1
$ws.Range("I31").Value = 1181.8206
$ws.Range("K31").Value = 1181.8206
$ws.Range("M31").Value = -886.8206

$ws.Range("H34").Value = 2046.1
$ws.Range("I34").Value = 1181.8206
$ws.Range("K34").Value = 1181.8206
$ws.Range("M34").Value = -979.8206

$ws.Range("H68").Value = 18970
$ws.Range("J68").Value = 18970
$ws.Range("L68").Value = 18970
$ws.Range("N68").Value = -20468

$ws.Range("H71").Value = 18970
$ws.Range("J71").Value = 18970
$ws.Range("L71").Value = 56910
$ws.Range("N71").Value = -64398

$ws.Range("H107").Value = 198.8125
$ws.Range("I107").Value = 170.91667
$ws.Range("J107").Value = 282.5
$ws.Range("K107").Value = 170.91667
$ws.Range("L107").Value = 282.5
$ws.Range("M107").Value = 1749.08333
$ws.Range("N107").Value = -4122.5

$ws.Range("H122").Value = 1972.4736
$ws.Range("I122").Value = 2299.7693
$ws.Range("J122").Value = 1263.3334
$ws.Range("K122").Value = 6899.3079
$ws.Range("L122").Value = 3790.0002
$ws.Range("M122").Value = -4449.3079
$ws.Range("N122").Value = -8690.0002

$ws.Range("H134").Value = 1646.4849
$ws.Range("I134").Value = 1120.4615
$ws.Range("K134").Value = 3361.3845
$ws.Range("M134").Value = -826.3844999999997

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 8.074074
$ws.Range("J12").Value = 2.6666667
$ws.Range("L12").Value = 8.000000099999999
$ws.Range("N12").Value = -354.0000001

$ws.Range("H23").Value = 443.35715
$ws.Range("I23").Value = 96.888885
$ws.Range("J23").Value = 607.4737
$ws.Range("K23").Value = 290.666655
$ws.Range("L23").Value = 1822.4211
$ws.Range("M23").Value = -55.66665499999999
$ws.Range("N23").Value = -2292.4211

$ws.Range("H121").Value = 2687.6206
$ws.Range("I121").Value = 4416.125
$ws.Range("J121").Value = 2029.1428
$ws.Range("K121").Value = 13248.375
$ws.Range("L121").Value = 6087.428400000001
$ws.Range("M121").Value = -11938.375
$ws.Range("N121").Value = -8707.428400000001

$ws.Range("H132").Value = 4420.231
$ws.Range("I132").Value = 2107.4119
$ws.Range("J132").Value = 8788.888999999999
$ws.Range("K132").Value = 18966.7071
$ws.Range("L132").Value = 79100.00099999999
$ws.Range("M132").Value = -16436.7071
$ws.Range("N132").Value = -84160.00099999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4678.579
$ws.Range("I70").Value = 4129.533
$ws.Range("K70").Value = 4129.533
$ws.Range("M70").Value = -3859.533

$ws.Range("H73").Value = 4678.579
$ws.Range("I73").Value = 4129.533
$ws.Range("K73").Value = 4129.533
$ws.Range("M73").Value = -3193.533

$ws.Range("H102").Value = 2546.6858
$ws.Range("I102").Value = 1522.4166
$ws.Range("J102").Value = 4781.4546
$ws.Range("K102").Value = 1522.4166
$ws.Range("L102").Value = 4781.4546
$ws.Range("M102").Value = 99.58339999999998
$ws.Range("N102").Value = -8025.4546

$ws.Range("H126").Value = 2550.8708
$ws.Range("I126").Value = 1948.2142
$ws.Range("J126").Value = 3047.1765
$ws.Range("K126").Value = 5844.642599999999
$ws.Range("L126").Value = 9141.529500000001
$ws.Range("M126").Value = -3374.642599999999
$ws.Range("N126").Value = -14081.5295

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1305.4
$ws.Range("I61").Value = 1381.75
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 1381.75
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -1179.75
$ws.Range("N61").Value = -1404

$ws.Range("H82").Value = 2270.6667
$ws.Range("I82").Value = 1658
$ws.Range("J82").Value = 3036.5
$ws.Range("K82").Value = 1658
$ws.Range("L82").Value = 3036.5
$ws.Range("M82").Value = -1297
$ws.Range("N82").Value = -3758.5

$ws.Range("H85").Value = 2270.6667
$ws.Range("I85").Value = 1658
$ws.Range("J85").Value = 3036.5
$ws.Range("K85").Value = 1658
$ws.Range("L85").Value = 3036.5
$ws.Range("M85").Value = -410
$ws.Range("N85").Value = -5532.5

$ws.Range("H113").Value = 1305.4
$ws.Range("I113").Value = 1381.75
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1381.75
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 788.25
$ws.Range("N113").Value = -5340

$ws.Range("H132").Value = 15505.8125
$ws.Range("I132").Value = 4749.125
$ws.Range("J132").Value = 26262.5
$ws.Range("K132").Value = 14247.375
$ws.Range("L132").Value = 78787.5
$ws.Range("M132").Value = -11717.375
$ws.Range("N132").Value = -83847.5

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
# Profit recalculated to 0/blank for this row; N134 drops out entirely.
$ws.Range("N134").ClearContents()
